{"js": "// fix: invisible navbar menu\n// The \"Home\" page-requirement heading and the paragraph describing it\n// (\"Frawlsalla wants this page to provide ...\") were accidentally left\n// highlighted in yellow. Remove the yellow highlight from both paragraphs\n// (including the paragraph mark) while leaving every other run property\n// (bold, bCs, etc.) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two paragraphs by their (unique) text content rather than a\n// fixed index, so the script is resilient to unrelated edits elsewhere in\n// the document.\nlet homeHeading = null;\nlet wantsParagraph = null;\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (text === \"Home\" && homeHeading === null) {\n    homeHeading = p;\n  } else if (text.indexOf(\"Frawlsalla wants this page\") === 0 && wantsParagraph === null) {\n    wantsParagraph = p;\n  }\n}\n\nif (!homeHeading || !wantsParagraph) {\n  throw new Error(\"Could not locate the 'Home' requirement paragraphs to fix the highlight.\");\n}\n\n// Setting font.highlightColor on the Paragraph proxy clears the highlight\n// from every run AND the paragraph mark's run properties (w:pPr/w:rPr),\n// which is exactly what the original fix did \u2014 drop <w:highlight .../>\n// everywhere in these two paragraphs while keeping bold/bCs intact.\nhomeHeading.font.highlightColor = null;\nwantsParagraph.font.highlightColor = null;\n\nawait context.sync();\n", "ps1": "# fix: invicible navbar menu\n#\n# The \"Home\" page-requirement heading and the paragraph describing it\n# (\"Frawlsalla wants this page to provide ...\") were accidentally left\n# highlighted in yellow, making the navbar-menu requirement text hard to\n# read. Remove the yellow highlight from both paragraphs (including the\n# paragraph mark) while leaving every other run property (bold, bCs, etc.)\n# untouched.\n\n$d = $word.ActiveDocument\n\n$homePara = $null\n$wantsPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n    if ($homePara -eq $null -and $text -eq \"Home\") {\n        $homePara = $p\n    } elseif ($wantsPara -eq $null -and $text.StartsWith(\"Frawlsalla wants this page\")) {\n        $wantsPara = $p\n    }\n}\n\n# wdNoHighlight = 0 \u2014 clearing Font.HighlightColorIndex on the paragraph's\n# full Range (which includes the trailing paragraph mark) drops\n# <w:highlight .../> from every run AND the paragraph mark's rPr, matching\n# the original fix.\nif ($homePara -ne $null) {\n    $homePara.Range.Font.HighlightColorIndex = 0\n}\nif ($wantsPara -ne $null) {\n    $wantsPara.Range.Font.HighlightColorIndex = 0\n}\n"}
